# Auto-generated edit script applying the diff to Sheets ALC, ARM, BSM, CRP, GSM, WVR
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 11111661
$ws.Range("I103").Value = 413.83334
$ws.Range("J103").Value = 16667285
$ws.Range("K103").Value = 1241.50002
$ws.Range("L103").Value = 50001855
$ws.Range("M103").Value = -655.5000199999999
$ws.Range("N103").Value = -50003027
# Row 120
$ws.Range("H120").Value = 36245
$ws.Range("J120").Value = 36245
$ws.Range("L120").Value = 36245
$ws.Range("N120").Value = -45921
# Row 132
$ws.Range("H132").Value = 3206.818
$ws.Range("I132").Value = 2593.7036
$ws.Range("J132").Value = 5965.8335
$ws.Range("K132").Value = 7781.110799999999
$ws.Range("L132").Value = 17897.5005
$ws.Range("M132").Value = -5251.110799999999
$ws.Range("N132").Value = -22957.5005

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 36029.74
$ws.Range("I32").Value = 8240.485000000001
$ws.Range("K32").Value = 8240.485000000001
$ws.Range("M32").Value = -7953.485000000001
# Row 55
$ws.Range("H55").Value = 14452
$ws.Range("J55").Value = 14452
$ws.Range("L55").Value = 14452
$ws.Range("N55").Value = -15082
# Row 61
$ws.Range("H61").Value = 3497.5
$ws.Range("I61").Value = 3497.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3497.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3285.5
$ws.Range("N61").ClearContents()
# Row 80
$ws.Range("H80").Value = 21261.25
$ws.Range("J80").Value = 21261.25
$ws.Range("L80").Value = 21261.25
$ws.Range("N80").Value = -23257.25
# Row 83
$ws.Range("H83").Value = 21261.25
$ws.Range("J83").Value = 21261.25
$ws.Range("L83").Value = 63783.75
$ws.Range("N83").Value = -73767.75
# Row 117
$ws.Range("H117").Value = 31667.2
$ws.Range("J117").Value = 31667.2
$ws.Range("L117").Value = 31667.2
$ws.Range("N117").Value = -40845.2
# Row 132
$ws.Range("H132").Value = 2300.617
$ws.Range("I132").Value = 1460.775
$ws.Range("K132").Value = 4382.325000000001
$ws.Range("M132").Value = -1852.325000000001
# Row 136
$ws.Range("H136").Value = 3497.5
$ws.Range("I136").Value = 3497.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10492.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7942.5
$ws.Range("N136").ClearContents()
# Row 140
$ws.Range("H140").Value = 73850.86
$ws.Range("J140").Value = 73850.86
$ws.Range("L140").Value = 73850.86
$ws.Range("N140").Value = -84210.86

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 24107.238
$ws.Range("J82").Value = 30814.666
$ws.Range("L82").Value = 30814.666
$ws.Range("N82").Value = -31580.666
# Row 85
$ws.Range("H85").Value = 24107.238
$ws.Range("J85").Value = 30814.666
$ws.Range("L85").Value = 30814.666
$ws.Range("N85").Value = -33466.666
# Row 94
$ws.Range("H94").Value = 1509
$ws.Range("I94").Value = 1509
$ws.Range("K94").Value = 1509
$ws.Range("M94").Value = -1058
# Row 105
$ws.Range("H105").Value = 2573.8386
$ws.Range("I105").Value = 1753.7084
$ws.Range("J105").Value = 5385.7144
$ws.Range("K105").Value = 1753.7084
$ws.Range("L105").Value = 5385.7144
$ws.Range("M105").Value = -6.708399999999983
$ws.Range("N105").Value = -8879.714400000001
# Row 134
$ws.Range("H134").Value = 1376.4166
$ws.Range("I134").Value = 1358.5
$ws.Range("J134").Value = 1466
$ws.Range("K134").Value = 4075.5
$ws.Range("L134").Value = 4398
$ws.Range("M134").Value = -1540.5
$ws.Range("N134").Value = -9468

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 16133.75
$ws.Range("J41").Value = 21345
$ws.Range("L41").Value = 21345
$ws.Range("N41").Value = -22201
# Row 50
$ws.Range("H50").Value = 8869.875
$ws.Range("I50").Value = 8083
$ws.Range("J50").Value = 8982.286
$ws.Range("K50").Value = 8083
$ws.Range("L50").Value = 8982.286
$ws.Range("M50").Value = -7458
$ws.Range("N50").Value = -10232.286
# Row 51
$ws.Range("H51").Value = 8148.4
$ws.Range("I51").Value = 3090
$ws.Range("J51").Value = 9413
$ws.Range("K51").Value = 3090
$ws.Range("L51").Value = 9413
$ws.Range("M51").Value = -2354
$ws.Range("N51").Value = -10885
# Row 60
$ws.Range("H60").Value = 22053.555
$ws.Range("J60").Value = 24423.625
$ws.Range("L60").Value = 24423.625
$ws.Range("N60").Value = -25445.625
# Row 61
$ws.Range("H61").Value = 8148.4
$ws.Range("I61").Value = 3090
$ws.Range("J61").Value = 9413
$ws.Range("K61").Value = 3090
$ws.Range("L61").Value = 9413
$ws.Range("M61").Value = -2742
$ws.Range("N61").Value = -10109
# Row 109
$ws.Range("H109").Value = 10957.143
$ws.Range("J109").Value = 10957.143
$ws.Range("L109").Value = 10957.143
$ws.Range("N109").Value = -13037.143
# Row 132
$ws.Range("H132").Value = 2761.1428
$ws.Range("I132").Value = 2882
$ws.Range("J132").Value = 1190
$ws.Range("K132").Value = 8646
$ws.Range("L132").Value = 3570
$ws.Range("M132").Value = -6116
$ws.Range("N132").Value = -8630
# Row 134
$ws.Range("H134").Value = 6628.893
$ws.Range("I134").Value = 7412.9165
$ws.Range("J134").Value = 1924.75
$ws.Range("K134").Value = 22238.7495
$ws.Range("L134").Value = 5774.25
$ws.Range("M134").Value = -19703.7495
$ws.Range("N134").Value = -10844.25

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3355.3076
$ws.Range("I102").Value = 1642.7142
$ws.Range("J102").Value = 5353.3335
$ws.Range("K102").Value = 1642.7142
$ws.Range("L102").Value = 5353.3335
$ws.Range("M102").Value = -20.71419999999989
$ws.Range("N102").Value = -8597.333500000001
# Row 107
$ws.Range("H107").Value = 586.44446
$ws.Range("I107").Value = 501.53845
$ws.Range("J107").Value = 665.2857
$ws.Range("K107").Value = 501.53845
$ws.Range("L107").Value = 665.2857
$ws.Range("M107").Value = 1418.46155
$ws.Range("N107").Value = -4505.2857
# Row 123
$ws.Range("H123").Value = 23267.334
$ws.Range("J123").Value = 23267.334
$ws.Range("L123").Value = 23267.334
$ws.Range("N123").Value = -28167.334
# Row 132
$ws.Range("H132").Value = 3291.842
$ws.Range("I132").Value = 2631.5
$ws.Range("J132").Value = 3596.6155
$ws.Range("K132").Value = 7894.5
$ws.Range("L132").Value = 10789.8465
$ws.Range("M132").Value = -5364.5
$ws.Range("N132").Value = -15849.8465

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 32700
$ws.Range("J109").Value = 32700
$ws.Range("L109").Value = 32700
$ws.Range("N109").Value = -35474
# Row 113
$ws.Range("H113").Value = 46293.637
$ws.Range("I113").Value = 62815.625
$ws.Range("J113").Value = 2235
$ws.Range("K113").Value = 188446.875
$ws.Range("L113").Value = 6705
$ws.Range("M113").Value = -186276.875
$ws.Range("N113").Value = -11045
# Row 132
$ws.Range("H132").Value = 2255
$ws.Range("I132").Value = 2341.3333
$ws.Range("J132").Value = 2047.8
$ws.Range("K132").Value = 7023.999899999999
$ws.Range("L132").Value = 6143.4
$ws.Range("M132").Value = -4493.999899999999
$ws.Range("N132").Value = -11203.4
# Row 136
$ws.Range("H136").Value = 870.4
$ws.Range("I136").Value = 711.34375
$ws.Range("K136").Value = 2134.03125
$ws.Range("M136").Value = 415.96875

Write-Host "Edit complete"